# Auto-generated edit script: updates cryptos price/volume table
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.588.58"
$ws.Range("E2").Value = "  +0.07%  "
$ws.Range("D3").Value = "2.469.43"
$ws.Range("E3").Value = "  +0.29%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "'318.61"
$ws.Range("E5").Value = "  +1.13%  "
$ws.Range("D6").Value = "'92.22"
$ws.Range("E6").Value = "  +1.12%  "
$ws.Range("D7").Value = "'0.551"
$ws.Range("E7").Value = "  +0.45%  "
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("D9").Value = "'0.512"
$ws.Range("E9").Value = "  +0.07%  "
$ws.Range("D10").Value = "'0.0864"
$ws.Range("E10").Value = "  +8.98%  "
$ws.Range("D11").Value = "'32.90"
$ws.Range("E11").Value = "  +1.21%  "
$ws.Range("E12").Value = "  -0.03%  "
$ws.Range("D13").Value = "2.848.08"
$ws.Range("E13").Value = "  +0.19%  "
$ws.Range("E14").Value = "  +0.09%  "
$ws.Range("D15").Value = "'15.45"
$ws.Range("E15").Value = "  -2.16%  "
$ws.Range("D16").Value = "2.467.88"
$ws.Range("E16").Value = "  -1.19%  "
$ws.Range("D17").Value = "'0.792"
$ws.Range("E17").Value = "  +2.50%  "
$ws.Range("D18").Value = "41.532.59"
$ws.Range("E18").Value = "  -0.10%  "
$ws.Range("D19").Value = "'6.44"
$ws.Range("E19").Value = "  -0.85%  "
$ws.Range("D20").Value = "0.0₃0942"
$ws.Range("E20").Value = "  +0.49%  "
$ws.Range("D21").Value = "'70.77"
$ws.Range("E21").Value = "  +0.00%  "
$ws.Range("D22").Value = "'11.24"
$ws.Range("E22").Value = "  +0.04%  "
$ws.Range("D23").Value = "'239.87"
$ws.Range("E23").Value = "  +1.25%  "
$ws.Range("E24").Value = "  +0.91%  "
$ws.Range("D25").Value = "'1.96"
$ws.Range("E25").Value = "  +2.78%  "
$ws.Range("E26").Value = "  +0.07%  "
$ws.Range("D27").Value = "'24.77"
$ws.Range("E27").Value = "  +2.14%  "
$ws.Range("E28").Value = "  -0.66%  "
$ws.Range("D29").Value = "'9.67"
$ws.Range("E29").Value = "  +0.49%  "
$ws.Range("D30").Value = "'36.29"
$ws.Range("E30").Value = "  +3.23%  "
$ws.Range("D31").Value = "'157.20"
$ws.Range("E31").Value = "  +0.43%  "
$ws.Range("D32").Value = "'5.45"
$ws.Range("E32").Value = "  +0.43%  "
$ws.Range("E33").Value = "  -0.16%  "
$ws.Range("B34").Value = "Hedera"
$ws.Range("C34").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D34").Value = "'0.0761"
$ws.Range("E34").Value = "  +0.66%  "
$ws.Range("B35").Value = "WEMIXToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D35").Value = "'2.56"
$ws.Range("E35").Value = "  -0.65%  "
$ws.Range("D36").Value = "'17.23"
$ws.Range("E36").Value = "  +0.33%  "
$ws.Range("D37").Value = "'1.85"
$ws.Range("E37").Value = "  +4.64%  "
$ws.Range("B38").Value = "LidoDAOToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D38").Value = "'2.89"
$ws.Range("E38").Value = "  +1.29%  "
$ws.Range("B39").Value = "Stellar"
$ws.Range("C39").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D39").Value = "'0.116"
$ws.Range("E39").Value = "  +1.40%  "
$ws.Range("E40").Value = "  +1.36%  "
$ws.Range("E41").Value = "  +5.07%  "
$ws.Range("D42").Value = "'3.98"
$ws.Range("E42").Value = "  +0.17%  "
$ws.Range("D43").Value = "1.990.42"
$ws.Range("E43").Value = "  +1.70%  "
$ws.Range("E44").Value = "  +0.25%  "
$ws.Range("D45").Value = "'18.80"
$ws.Range("E45").Value = "  +0.41%  "
$ws.Range("D46").Value = "'2.95"
$ws.Range("E46").Value = "  +1.46%  "
$ws.Range("D47").Value = "'9.40"
$ws.Range("E47").Value = "  +4.74%  "
$ws.Range("D48").Value = "2.703.76"
$ws.Range("E48").Value = "  +0.07%  "
$ws.Range("D49").Value = "'97.45"
$ws.Range("E49").Value = "  +1.07%  "
$ws.Range("D50").Value = "'75.77"
$ws.Range("E50").Value = "  +5.97%  "
$ws.Range("D51").Value = "'66.74"
$ws.Range("E51").Value = "  -0.20%  "
